$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connectivity")
Write-Host $ws.Name
